# RemoveIncomeAndExpense.xlsx edit
# - fix typo in the shared string used by D4/E4
# - move the active selection from G2 to E9
# - nudge the A:B and D:E column widths
# - restore the workbook window geometry (best effort; engine may not
#   persist window placement, but set it for correctness anyway)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the typo: "โครงการะหายไป" -> "โครงการหายไป" (stray "ะ" removed) ---
$fixedText = "ข้อมูลรายรับรายจ่ายของโครงการหายไป"
$ws.Range("D4").Value = $fixedText
$ws.Range("E4").Value = $fixedText

# --- Column widths: A:B narrower, D:E wider ---
$ws.Columns.Item(1).ColumnWidth = 7.8
$ws.Columns.Item(2).ColumnWidth = 7.8
$ws.Columns.Item(4).ColumnWidth = 28.8
$ws.Columns.Item(5).ColumnWidth = 28.8

# --- Move the selection/active cell to E9 ---
$ws.Range("E9").Select() | Out-Null

# --- Restore workbook window placement/size ---
$win = $excel.ActiveWindow
$win.Left = 150
$win.Top = 2610
$win.Width = 12585
$win.Height = 11235
